$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.721.75'
$ws.Range("E2").Value = '  +6.92%  '

$ws.Range("D3").Value = '2.037.58'
$ws.Range("E3").Value = '  +8.06%  '

$ws.Range("E4").Value = '  +0.00%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '254.60'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.92%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.702'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.10%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.07%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '47.98'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +12.46%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.386'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +9.36%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '58.64'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +6.43%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0775'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +4.68%  '

$ws.Range("E12").Value = '  +2.83%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '15.58'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +13.23%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.852'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +9.43%  '

$ws.Range("D15").Value = '2.310.14'
$ws.Range("E15").Value = '  +7.05%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.25'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +5.61%  '

$ws.Range("D17").Value = '2.023.02'
$ws.Range("E17").Value = '  +7.35%  '

$ws.Range("D18").Value = '37.749.69'
$ws.Range("E18").Value = '  +7.07%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '75.83'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +3.61%  '

$ws.Range("E20").Value = '  +5.86%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.92'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +9.18%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '257.28'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +5.40%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.34'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +4.07%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.998'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.19%  '

$ws.Range("E25").Value = '  -5.49%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '170.73'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.22%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.17'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.44%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.03'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +6.45%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '20.44'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +12.08%  '

$ws.Range("E30").Value = '  +2.79%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.63'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +8.53%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.0623'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +5.13%  '

$ws.Range("B33").Value = 'Gas'
$ws.Range("C33").Value = 'https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '21.35'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +61.54%  '

$ws.Range("B34").Value = 'Kaspa'
$ws.Range("C34").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0923'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +28.71%  '

$ws.Range("B35").Value = 'InternetComputer(DFINITY)'
$ws.Range("C35").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.40'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +5.85%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.89'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.15%  '

$ws.Range("E37").Value = '  +0.05%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.913'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +8.19%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.46'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.63%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.18'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +12.91%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '104.59'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +7.20%  '

$ws.Range("E42").Value = '  +4.57%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '17.71'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +3.67%  '

$ws.Range("B44").Value = 'ARBITRUM'
$ws.Range("C44").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.14'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +7.18%  '

$ws.Range("B45").Value = 'HuobiToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.90'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +20.73%  '

$ws.Range("D46").Value = '1.389.94'
$ws.Range("E46").Value = '  +4.81%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.46'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +4.78%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0857'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +6.09%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.87'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +4.68%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '3.98'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +20.46%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '6.60'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +5.36%  '
